# Apply corrections to the Diebold-Mariano summary table (dm_summary_General.xlsx)
# Changes:
#  - B4 (Sieve Bootstrap, Comparaciones_Significativas): "4/10" -> "5/10"
#  - C4 (Sieve Bootstrap, Proporcion_Sig): 102.4 -> 128
#  - C6 (AV-MCPS, Proporcion_Sig): 25.6 -> 0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "5/10"
$ws.Range("C4").Value = 128
$ws.Range("C6").Value = 0
